# Added code for New Task functionality
# - Adds "FeedbackForms" and "Tasks" worksheets (with sample data) to the
#   freeCrmTestData workbook, mirroring the existing test-data sheets.

$wb = $excel.ActiveWorkbook

# --- MultiDocs: row 1 gets selected (no longer the active tab once the new
#     sheets are appended) -----------------------------------------------
$wsMultiDocs = $wb.Worksheets.Item("MultiDocs")
$wsMultiDocs.Rows.Item(1).Select() | Out-Null

# --- Add the two new worksheets, in order, after MultiDocs ---------------
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForms = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLast)
$wsForms.Name = "FeedbackForms"

$wsTasks = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsForms)
$wsTasks.Name = "Tasks"

# =====================================================================
# FeedbackForms sheet
# =====================================================================

# Header row
$wsForms.Range("A1").Value = "title"
$wsForms.Range("B1").Value = "pages"
$wsForms.Range("C1").Value = "reportEmail"
$wsForms.Range("D1").Value = "Description"
$wsForms.Range("E1").Value = "welcomeMessage"
$wsForms.Range("F1").Value = "confirmationMessage"
$wsForms.Range("A1:F1").Interior.Color = 65535

# Row 2
$wsForms.Range("A2").Value = "Form title 1"
$wsForms.Range("B2").Value = "'10"
$wsForms.Range("C2").Value = "abc@gmail.com"
$wsForms.Range("D2").Value = "Form 1 description"
$wsForms.Range("E2").Value = "Form 1 welcome message"
$wsForms.Range("F2").Value = "Form1 confirmation message"

# Row 3
$wsForms.Range("A3").Value = "Form title 2"
$wsForms.Range("B3").Value = "'15"
$wsForms.Range("C3").Value = "def@gmail.com"
$wsForms.Range("D3").Value = "Form 2 description"
$wsForms.Range("E3").Value = "Form 2 welcome message"
$wsForms.Range("F3").Value = "Form2 confirmation message"

# Hyperlinks for the report e-mail addresses
$wsForms.Hyperlinks.Add($wsForms.Range("C2"), "mailto:abc@gmail.com") | Out-Null
$wsForms.Hyperlinks.Add($wsForms.Range("C3"), "mailto:def@gmail.com") | Out-Null

$wsForms.Range("F3").Select() | Out-Null

# =====================================================================
# Tasks sheet
# =====================================================================

# Header row
$wsTasks.Range("A1").Value = "title"
$wsTasks.Range("B1").Value = "autoextend"
$wsTasks.Range("C1").Value = "status"
$wsTasks.Range("D1").Value = "completionpercentage"
$wsTasks.Range("E1").Value = "type"
$wsTasks.Range("F1").Value = "priority"
$wsTasks.Range("G1").Value = "deal"
$wsTasks.Range("H1").Value = "case"
$wsTasks.Range("I1").Value = "tags"
$wsTasks.Range("J1").Value = "description"
$wsTasks.Range("K1").Value = "contact"
$wsTasks.Range("L1").Value = "company"
$wsTasks.Range("M1").Value = "identifier"
$wsTasks.Range("A1:M1").Interior.Color = 65535

# Data rows, filled column by column (column A rows 2-3, then column B
# rows 2-3, ...) to match the shared-string insertion order of the source
# workbook.
$wsTasks.Range("A2").Value = "Task1"
$wsTasks.Range("A3").Value = "Task2"

$wsTasks.Range("B2").Value = "Extend deadline by 1 day"
$wsTasks.Range("B3").Value = "Extend deadline by 30 days"

$wsTasks.Range("C2").Value = "Open"
$wsTasks.Range("C3").Value = "Complete"

$wsTasks.Range("D2").Value = "'50"
$wsTasks.Range("D3").Value = "'75"

$wsTasks.Range("E2").Value = "Call"
$wsTasks.Range("E3").Value = "Training"

$wsTasks.Range("F2").Value = "High"
$wsTasks.Range("F3").Value = "Normal"

$wsTasks.Range("G2").Value = "deal1"
$wsTasks.Range("G3").Value = "deal2"

$wsTasks.Range("H2").Value = "case1"
$wsTasks.Range("H3").Value = "case2"

$wsTasks.Range("I2").Value = "tag1"
$wsTasks.Range("I3").Value = "tag2"

$wsTasks.Range("J2").Value = "tagdesc1"
$wsTasks.Range("J3").Value = "tagdesc2"

$wsTasks.Range("K2").Value = "contact1"
$wsTasks.Range("K3").Value = "contact2"

$wsTasks.Range("L2").Value = "company1"
$wsTasks.Range("L3").Value = "company2"

$wsTasks.Range("M2").Value = "identifier1"
$wsTasks.Range("M3").Value = "identifier2"

$wsTasks.Columns.Item(2).ColumnWidth = 22.333333333333332
$wsTasks.Columns.Item(4).ColumnWidth = 18.666666666666668

$wsTasks.Range("I19").Select() | Out-Null
